$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.029.85'
$ws.Range('E2').Value = '  -1.55%  '
$ws.Range('D3').Value = '2.416.70'
$ws.Range('E3').Value = '  -2.35%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.86'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.79%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '88.65'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.51%  '
$ws.Range('E7').Value = '  -2.73%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.496'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.96%  '
$ws.Range('E10').Value = '  -2.16%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '31.51'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.92%  '
$ws.Range('E12').Value = '  -1.62%  '
$ws.Range('D13').Value = '2.787.15'
$ws.Range('E13').Value = '  -2.40%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.78'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.78%  '
$ws.Range('E15').Value = '  -1.12%  '
$ws.Range('D16').Value = '2.413.19'
$ws.Range('E16').Value = '  -2.30%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.771'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.33%  '
$ws.Range('D18').Value = '40.843.07'
$ws.Range('E18').Value = '  -1.90%  '
$ws.Range('D19').Value = '0.0₃0921'
$ws.Range('E19').Value = '  -3.24%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.24'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.68%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '70.92'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.50%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.88'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.85%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '238.62'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.40%  '
$ws.Range('E24').Value = '  -3.11%  '
$ws.Range('E25').Value = '  +0.20%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.85'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.36%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.04'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.18%  '
$ws.Range('E28').Value = '  -2.57%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.52'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.42%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '34.22'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.40%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '157.69'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.96%  '
$ws.Range('E32').Value = '  -4.61%  '
$ws.Range('E33').Value = '  +0.05%  '
$ws.Range('E34').Value = '  -4.02%  '
$ws.Range('E35').Value = '  -4.82%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.89'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.34%  '
$ws.Range('B37').Value = 'Stellar'
$ws.Range('C37').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.114'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.23%  '
$ws.Range('B38').Value = 'Celestia'
$ws.Range('C38').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '16.32'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.77%  '
$ws.Range('E39').Value = '  -6.13%  '
$ws.Range('E40').Value = '  -3.56%  '
$ws.Range('E41').Value = '  -3.26%  '
$ws.Range('E42').Value = '  -6.71%  '
$ws.Range('D43').Value = '1.987.46'
$ws.Range('E43').Value = '  -0.31%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0274'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.15%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.18'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.83%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.86'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.56%  '
$ws.Range('E47').Value = '  +0.47%  '
$ws.Range('D48').Value = '2.647.52'
$ws.Range('E48').Value = '  -2.37%  '
$ws.Range('B49').Value = 'BitcoinSV'
$ws.Range('C49').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '73.55'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.77%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '94.06'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.37%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '51.36'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.77%  '
